$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1440.690759875038
$ws.Range("E2").Value = 1707.8679159935714
$ws.Range("G2").Value = 955.5854869749402
$ws.Range("J2").Value = 488.54336074569926
$ws.Range("C3").Value = 1102.6446499341537
$ws.Range("E3").Value = 1822.2997128583345
$ws.Range("G3").Value = 1194.1015187255632
$ws.Range("J3").Value = 519.3330594676355
$ws.Range("C4").Value = 1112.3315134753032
$ws.Range("E4").Value = 2035.6893118204573
$ws.Range("G4").Value = 799.1835329027557
$ws.Range("J4").Value = 811.4804669414951
$ws.Range("C5").Value = 1100.2900609715684
$ws.Range("E5").Value = 1722.9788873957802
$ws.Range("G5").Value = 895.2177544738933
$ws.Range("J5").Value = 739.349160139947
$ws.Range("C6").Value = 546.5450380717565
$ws.Range("E6").Value = 1931.5995437187648
$ws.Range("G6").Value = 741.3650300668102
$ws.Range("J6").Value = 152.5266386930083
$ws.Range("C7").Value = 801.8704103813825
$ws.Range("E7").Value = 1713.9983884607057
$ws.Range("G7").Value = 708.4912224940716
$ws.Range("J7").Value = 516.8625075072472
$ws.Range("C8").Value = 929.5909822242468
$ws.Range("E8").Value = 1884.9234449812907
$ws.Range("G8").Value = 701.1341987630879
$ws.Range("J8").Value = 698.1750159023352
$ws.Range("C9").Value = 969.4593980148438
$ws.Range("E9").Value = 1869.9693203495092
$ws.Range("G9").Value = 743.9308036328167
$ws.Range("J9").Value = 818.8938606109518
$ws.Range("C10").Value = 1059.773744490573
$ws.Range("E10").Value = 2510.7021607639144
$ws.Range("G10").Value = 858.83915007571
$ws.Range("J10").Value = 778.1921690866949
$ws.Range("C11").Value = 1142.2503670735578
$ws.Range("E11").Value = 1705.1364138896488
$ws.Range("G11").Value = 753.2402760522962
$ws.Range("J11").Value = 707.7514220216078
